$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.693.43"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "3.867.94"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.58"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.82"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.70%  "
$ws.Range("D7").Value = "3.869.25"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("E10").Value = "  +3.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.50"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.42%  "
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000288"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +16.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.23"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("D15").Value = "4.511.38"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "3.857.49"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "68.697.35"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.40"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.08"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "472.97"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.733"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("E24").Value = "  +2.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.91"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("E26").Value = "  +3.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.25"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.50"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.49%  "
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").Value = "4.017.49"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.78"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.39"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("D36").Value = "3.830.32"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.96"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +21.71%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.01"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.02"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.10%  "
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  +3.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000303"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +14.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.83"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.59%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "420.08"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.67"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.76"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("E51").Value = "  +1.68%  "
